$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The five "Knärot" observation rows (17-21) got re-keyed: each row's Id (A)
# and coordinate pair (Q,R) now carries the data that used to live one slot
# over in the cycle 17 -> 18 -> 20 -> 21 -> 19 -> 17. Row 19's Enhet/Ålder-
# Stadium/Kön/Metod/Bestämningsmetod annotations (J/K/L/N/AF, with
# K="blomning") move down to row 20.

# Row 17 <- old row 18
$ws.Range("A17").Value = 111821927
$ws.Range("Q17").Value = 550819.8901872271
$ws.Range("R17").Value = 6681733.007140613

# Row 18 <- old row 20
$ws.Range("A18").Value = 111821928
$ws.Range("Q18").Value = 550825.9503372401
$ws.Range("R18").Value = 6681726.144349095

# Row 19 <- old row 17 (loses the J/K/L/N/AF annotation cells)
$ws.Range("A19").Value = 111821926
$ws.Range("Q19").Value = 550846.2444635418
$ws.Range("R19").Value = 6681625.195240833
$ws.Range("J19:N19").ClearContents()
$ws.Range("AF19").ClearContents()

# Row 20 <- old row 21 (gains the J/K/L/N/AF annotation cells)
$ws.Range("A20").Value = 111821923
$ws.Range("Q20").Value = 550701.1291094749
$ws.Range("R20").Value = 6681909.496304798
$ws.Range("J20").Value = ""
$ws.Range("K20").Value = "blomning"
$ws.Range("L20").Value = ""
$ws.Range("N20").Value = ""
$ws.Range("AF20").Value = ""

# Row 21 <- old row 19 (keeps its own J/K/L/N/AF annotation cells as-is)
$ws.Range("A21").Value = 111821924
$ws.Range("Q21").Value = 550675.3931295178
$ws.Range("R21").Value = 6681937.422269406
